$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BF column holds a "game date" stamp for each team's row. Because of how
# the NBA stats site displayed dates, the value was written in a malformed,
# off-by-a-day form like "4-29-2012-13" instead of the correct ISO-style
# "2013-04-29". Fix every row that still has the bad value.
#
# NOTE: Assigning a date-shaped string straight to .Value/.Value2/.Formula
# makes Excel "helpfully" reinterpret it as an actual date serial number
# (e.g. 41393) instead of literal text. To keep it as plain text (matching
# how the column was already stored), write it as a text-literal formula
# ("=""2013-04-29""") and then collapse the formula down to a static value
# with Copy / PasteSpecial values-only - this avoids Excel's automatic
# date-detection while leaving every other formatting/style untouched.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)  # column BF
    if ($cell.Value2 -eq "4-29-2012-13") {
        $cell.Formula = '="2013-04-29"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = 0
